$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.317248821258545
$ws.Range("B1").Value = 6.210578441619873
$ws.Range("C1").Value = 5.124565601348877
$ws.Range("D1").Value = 5.956264972686768
$ws.Range("E1").Value = 4.773531913757324
